$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 223, shifting rows 223:297 down to 224:298
$ws.Rows("223:223").Insert()

# Populate the newly inserted row 223 with the new data
$ws.Cells.Item(223, 1).Value = 3
$ws.Cells.Item(223, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(223, 3).Value = "Coquimbo"
$ws.Cells.Item(223, 4).Value = 44627
$ws.Cells.Item(223, 5).Value = 5
$ws.Cells.Item(223, 6).Value = 100112012
$ws.Cells.Item(223, 7).Value = "Espinaca"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 126
$ws.Cells.Item(223, 11).Value = 4500
$ws.Cells.Item(223, 12).Value = 5000
$ws.Cells.Item(223, 13).Value = 4722
$ws.Cells.Item(223, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(223, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(223, 16).Value = 1574
$ws.Cells.Item(223, 17).Value = 3
$ws.Cells.Item(223, 18).Value = "Hortaliza"
